$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Wnt4'
$ws.Cells.Item(2, 3).Value = 'Fzd6'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.285322
$ws.Cells.Item(2, 8).Value = 9.855966
$ws.Cells.Item(2, 9).Value = 0.4533344065718998
$ws.Cells.Item(2, 10).Value = 0.4533344065718997
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 13.21223933333333
$ws.Cells.Item(2, 14).Value = 39.636718
$ws.Cells.Item(2, 15).Value = 0.8149747101495924
$ws.Cells.Item(2, 16).Value = 0.8149747101495927
$ws.Cells.Item(2, 17).Value = 43.40646055106534
$ws.Cells.Item(2, 18).Value = 390.658144959588
$ws.Cells.Item(2, 19).Value = 0.3694560765967715
$ws.Cells.Item(2, 20).Value = 0.3694560765967715

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Wnt4'
$ws.Cells.Item(3, 3).Value = 'Fzd6'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.285322
$ws.Cells.Item(3, 8).Value = 9.855966
$ws.Cells.Item(3, 9).Value = 0.4533344065718998
$ws.Cells.Item(3, 10).Value = 0.4533344065718997
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.662736333333334
$ws.Cells.Item(3, 14).Value = 7.988209000000001
$ws.Cells.Item(3, 15).Value = 0.1642464018940561
$ws.Cells.Item(3, 16).Value = 0.1642464018940561
$ws.Cells.Item(3, 17).Value = 8.747946256099334
$ws.Cells.Item(3, 18).Value = 78.73151630489401
$ws.Cells.Item(3, 19).Value = 0.07445854513421166
$ws.Cells.Item(3, 20).Value = 0.07445854513421167

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Wnt4'
$ws.Cells.Item(4, 3).Value = 'Fzd6'
$ws.Cells.Item(4, 4).Value = 'M2'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.285322
$ws.Cells.Item(4, 8).Value = 9.855966
$ws.Cells.Item(4, 9).Value = 0.4533344065718998
$ws.Cells.Item(4, 10).Value = 0.4533344065718997
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.009795
$ws.Cells.Item(4, 14).Value = 0.029385
$ws.Cells.Item(4, 15).Value = 0.0006041880626379251
$ws.Cells.Item(4, 16).Value = 0.0006041880626379253
$ws.Cells.Item(4, 17).Value = 0.03217972899
$ws.Cells.Item(4, 18).Value = 0.28961756091
$ws.Cells.Item(4, 19).Value = 0.0002738992368337896
$ws.Cells.Item(4, 20).Value = 0.0002738992368337896

# Row 5
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Wnt4'
$ws.Cells.Item(5, 3).Value = 'Fzd6'
$ws.Cells.Item(5, 4).Value = 'sCs'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.285322
$ws.Cells.Item(5, 8).Value = 9.855966
$ws.Cells.Item(5, 9).Value = 0.4533344065718998
$ws.Cells.Item(5, 10).Value = 0.4533344065718997
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.327069
$ws.Cells.Item(5, 14).Value = 0.9812069999999999
$ws.Cells.Item(5, 15).Value = 0.02017469989371348
$ws.Cells.Item(5, 16).Value = 0.02017469989371348
$ws.Cells.Item(5, 17).Value = 1.074526981218
$ws.Cells.Item(5, 18).Value = 9.670742830962
$ws.Cells.Item(5, 19).Value = 0.009145885604082769
$ws.Cells.Item(5, 20).Value = 0.009145885604082769

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Wnt4'
$ws.Cells.Item(6, 3).Value = 'Fzd6'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.533623
$ws.Cells.Item(6, 8).Value = 4.600869
$ws.Cells.Item(6, 9).Value = 0.2116212878402838
$ws.Cells.Item(6, 10).Value = 0.2116212878402837
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 13.21223933333333
$ws.Cells.Item(6, 14).Value = 39.636718
$ws.Cells.Item(6, 15).Value = 0.8149747101495924
$ws.Cells.Item(6, 16).Value = 0.8149747101495927
$ws.Cells.Item(6, 17).Value = 20.26259412310467
$ws.Cells.Item(6, 18).Value = 182.363347107942
$ws.Cells.Item(6, 19).Value = 0.1724659977191187
$ws.Cells.Item(6, 20).Value = 0.1724659977191187

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Wnt4'
$ws.Cells.Item(7, 3).Value = 'Fzd6'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.533623
$ws.Cells.Item(7, 8).Value = 4.600869
$ws.Cells.Item(7, 9).Value = 0.2116212878402838
$ws.Cells.Item(7, 10).Value = 0.2116212878402837
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.662736333333334
$ws.Cells.Item(7, 14).Value = 7.988209000000001
$ws.Cells.Item(7, 15).Value = 0.1642464018940561
$ws.Cells.Item(7, 16).Value = 0.1642464018940561
$ws.Cells.Item(7, 17).Value = 4.083633683735668
$ws.Cells.Item(7, 18).Value = 36.75270315362101
$ws.Cells.Item(7, 19).Value = 0.03475803509195297
$ws.Cells.Item(7, 20).Value = 0.03475803509195297

# Row 8
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Wnt4'
$ws.Cells.Item(8, 3).Value = 'Fzd6'
$ws.Cells.Item(8, 4).Value = 'M2'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.533623
$ws.Cells.Item(8, 8).Value = 4.600869
$ws.Cells.Item(8, 9).Value = 0.2116212878402838
$ws.Cells.Item(8, 10).Value = 0.2116212878402837
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009795
$ws.Cells.Item(8, 14).Value = 0.029385
$ws.Cells.Item(8, 15).Value = 0.0006041880626379251
$ws.Cells.Item(8, 16).Value = 0.0006041880626379253
$ws.Cells.Item(8, 17).Value = 0.015021837285
$ws.Cells.Item(8, 18).Value = 0.135196535565
$ws.Cells.Item(8, 19).Value = 0.0001278590559131637
$ws.Cells.Item(8, 20).Value = 0.0001278590559131637

# Row 9
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Wnt4'
$ws.Cells.Item(9, 3).Value = 'Fzd6'
$ws.Cells.Item(9, 4).Value = 'sCs'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.533623
$ws.Cells.Item(9, 8).Value = 4.600869
$ws.Cells.Item(9, 9).Value = 0.2116212878402838
$ws.Cells.Item(9, 10).Value = 0.2116212878402837
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.327069
$ws.Cells.Item(9, 14).Value = 0.9812069999999999
$ws.Cells.Item(9, 15).Value = 0.02017469989371348
$ws.Cells.Item(9, 16).Value = 0.02017469989371348
$ws.Cells.Item(9, 17).Value = 0.501600540987
$ws.Cells.Item(9, 18).Value = 4.514404868883
$ws.Cells.Item(9, 19).Value = 0.004269395973298882
$ws.Cells.Item(9, 20).Value = 0.004269395973298882

# Row 10
$ws.Cells.Item(10, 1).Value = 'M2'
$ws.Cells.Item(10, 2).Value = 'Wnt4'
$ws.Cells.Item(10, 3).Value = 'Fzd6'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.389186666666667
$ws.Cells.Item(10, 8).Value = 4.16756
$ws.Cells.Item(10, 9).Value = 0.1916908336993843
$ws.Cells.Item(10, 10).Value = 0.1916908336993843
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 13.21223933333333
$ws.Cells.Item(10, 14).Value = 39.636718
$ws.Cells.Item(10, 15).Value = 0.8149747101495924
$ws.Cells.Item(10, 16).Value = 0.8149747101495927
$ws.Cells.Item(10, 17).Value = 18.35426671867556
$ws.Cells.Item(10, 18).Value = 165.18840046808
$ws.Cells.Item(10, 19).Value = 0.1562231816324895
$ws.Cells.Item(10, 20).Value = 0.1562231816324895

# Row 11
$ws.Cells.Item(11, 1).Value = 'M2'
$ws.Cells.Item(11, 2).Value = 'Wnt4'
$ws.Cells.Item(11, 3).Value = 'Fzd6'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.389186666666667
$ws.Cells.Item(11, 8).Value = 4.16756
$ws.Cells.Item(11, 9).Value = 0.1916908336993843
$ws.Cells.Item(11, 10).Value = 0.1916908336993843
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.662736333333334
$ws.Cells.Item(11, 14).Value = 7.988209000000001
$ws.Cells.Item(11, 15).Value = 0.1642464018940561
$ws.Cells.Item(11, 16).Value = 0.1642464018940561
$ws.Cells.Item(11, 17).Value = 3.699037811115556
$ws.Cells.Item(11, 18).Value = 33.29134030004001
$ws.Cells.Item(11, 19).Value = 0.03148452971119575
$ws.Cells.Item(11, 20).Value = 0.03148452971119576

# Row 12
$ws.Cells.Item(12, 1).Value = 'M2'
$ws.Cells.Item(12, 2).Value = 'Wnt4'
$ws.Cells.Item(12, 3).Value = 'Fzd6'
$ws.Cells.Item(12, 4).Value = 'M2'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.389186666666667
$ws.Cells.Item(12, 8).Value = 4.16756
$ws.Cells.Item(12, 9).Value = 0.1916908336993843
$ws.Cells.Item(12, 10).Value = 0.1916908336993843
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.009795
$ws.Cells.Item(12, 14).Value = 0.029385
$ws.Cells.Item(12, 15).Value = 0.0006041880626379251
$ws.Cells.Item(12, 16).Value = 0.0006041880626379253
$ws.Cells.Item(12, 17).Value = 0.0136070834
$ws.Cells.Item(12, 18).Value = 0.1224637506
$ws.Cells.Item(12, 19).Value = 0.0001158173134382797
$ws.Cells.Item(12, 20).Value = 0.0001158173134382797

# Row 13
$ws.Cells.Item(13, 1).Value = 'M2'
$ws.Cells.Item(13, 2).Value = 'Wnt4'
$ws.Cells.Item(13, 3).Value = 'Fzd6'
$ws.Cells.Item(13, 4).Value = 'sCs'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.389186666666667
$ws.Cells.Item(13, 8).Value = 4.16756
$ws.Cells.Item(13, 9).Value = 0.1916908336993843
$ws.Cells.Item(13, 10).Value = 0.1916908336993843
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.327069
$ws.Cells.Item(13, 14).Value = 0.9812069999999999
$ws.Cells.Item(13, 15).Value = 0.02017469989371348
$ws.Cells.Item(13, 16).Value = 0.02017469989371348
$ws.Cells.Item(13, 17).Value = 0.45435989388
$ws.Cells.Item(13, 18).Value = 4.089239044919999
$ws.Cells.Item(13, 19).Value = 0.003867305042260817
$ws.Cells.Item(13, 20).Value = 0.003867305042260818

# Row 14
$ws.Cells.Item(14, 1).Value = 'sCs'
$ws.Cells.Item(14, 2).Value = 'Wnt4'
$ws.Cells.Item(14, 3).Value = 'Fzd6'
$ws.Cells.Item(14, 4).Value = 'ECs'
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.038885
$ws.Cells.Item(14, 8).Value = 3.116655
$ws.Cells.Item(14, 9).Value = 0.1433534718884322
$ws.Cells.Item(14, 10).Value = 0.1433534718884322
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 13.21223933333333
$ws.Cells.Item(14, 14).Value = 39.636718
$ws.Cells.Item(14, 15).Value = 0.8149747101495924
$ws.Cells.Item(14, 16).Value = 0.8149747101495927
$ws.Cells.Item(14, 17).Value = 13.72599725981
$ws.Cells.Item(14, 18).Value = 123.53397533829
$ws.Cells.Item(14, 19).Value = 0.1168294542012128
$ws.Cells.Item(14, 20).Value = 0.1168294542012128

# Row 15
$ws.Cells.Item(15, 1).Value = 'sCs'
$ws.Cells.Item(15, 2).Value = 'Wnt4'
$ws.Cells.Item(15, 3).Value = 'Fzd6'
$ws.Cells.Item(15, 4).Value = 'FAPs'
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.038885
$ws.Cells.Item(15, 8).Value = 3.116655
$ws.Cells.Item(15, 9).Value = 0.1433534718884322
$ws.Cells.Item(15, 10).Value = 0.1433534718884322
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.662736333333334
$ws.Cells.Item(15, 14).Value = 7.988209000000001
$ws.Cells.Item(15, 15).Value = 0.1642464018940561
$ws.Cells.Item(15, 16).Value = 0.1642464018940561
$ws.Cells.Item(15, 17).Value = 2.766276835655
$ws.Cells.Item(15, 18).Value = 24.896491520895
$ws.Cells.Item(15, 19).Value = 0.02354529195669571
$ws.Cells.Item(15, 20).Value = 0.02354529195669572

# Row 16
$ws.Cells.Item(16, 1).Value = 'sCs'
$ws.Cells.Item(16, 2).Value = 'Wnt4'
$ws.Cells.Item(16, 3).Value = 'Fzd6'
$ws.Cells.Item(16, 4).Value = 'M2'
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.038885
$ws.Cells.Item(16, 8).Value = 3.116655
$ws.Cells.Item(16, 9).Value = 0.1433534718884322
$ws.Cells.Item(16, 10).Value = 0.1433534718884322
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.009795
$ws.Cells.Item(16, 14).Value = 0.029385
$ws.Cells.Item(16, 15).Value = 0.0006041880626379251
$ws.Cells.Item(16, 16).Value = 0.0006041880626379253
$ws.Cells.Item(16, 17).Value = 0.010175878575
$ws.Cells.Item(16, 18).Value = 0.09158290717499999
$ws.Cells.Item(16, 19).Value = 0.00008661245645269212
$ws.Cells.Item(16, 20).Value = 0.00008661245645269215

# Row 17
$ws.Cells.Item(17, 1).Value = 'sCs'
$ws.Cells.Item(17, 2).Value = 'Wnt4'
$ws.Cells.Item(17, 3).Value = 'Fzd6'
$ws.Cells.Item(17, 4).Value = 'sCs'
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.038885
$ws.Cells.Item(17, 8).Value = 3.116655
$ws.Cells.Item(17, 9).Value = 0.1433534718884322
$ws.Cells.Item(17, 10).Value = 0.1433534718884322
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.327069
$ws.Cells.Item(17, 14).Value = 0.9812069999999999
$ws.Cells.Item(17, 15).Value = 0.02017469989371348
$ws.Cells.Item(17, 16).Value = 0.02017469989371348
$ws.Cells.Item(17, 17).Value = 0.339787078065
$ws.Cells.Item(17, 18).Value = 3.058083702584999
$ws.Cells.Item(17, 19).Value = 0.002892113274071012
$ws.Cells.Item(17, 20).Value = 0.002892113274071012
